# Update the cryptos worksheet with the latest scraped market data.
# Values that look numeric (e.g. "1.00", "33.60") are written with a
# leading apostrophe so Excel keeps them as text instead of silently
# converting them to numbers and dropping the formatted trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'87.509.78"
$ws.Range("E2").Value = "  +4.06%  "
$ws.Range("D3").Value = "'3.267.67"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'212.37"
$ws.Range("E5").Value = "  -2.91%  "
$ws.Range("D6").Value = "'625.06"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("D7").Value = "'0.366"
$ws.Range("E7").Value = "  +13.67%  "
$ws.Range("D8").Value = "'0.684"
$ws.Range("E8").Value = "  +15.97%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'3.267.56"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").Value = "'0.572"
$ws.Range("E11").Value = "  -4.27%  "
$ws.Range("D12").Value = "'0.179"
$ws.Range("E12").Value = "  +7.54%  "
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  -8.45%  "
$ws.Range("D14").Value = "'3.857.97"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "'33.60"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "'86.935.69"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("D18").Value = "'3.251.15"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "'3.10"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").Value = "'13.99"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("D21").Value = "'432.67"
$ws.Range("E21").Value = "  -3.50%  "
$ws.Range("D22").Value = "'8.85"
$ws.Range("E22").Value = "  -3.27%  "
$ws.Range("D23").Value = "'5.31"
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("D24").Value = "'7.28"
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("D25").Value = "'12.44"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  -4.55%  "
$ws.Range("D27").Value = "'3.417.06"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").Value = "'75.86"
$ws.Range("E28").Value = "  -2.85%  "
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +12.13%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'8.70"
$ws.Range("E33").Value = "  -5.33%  "
$ws.Range("D34").Value = "'542.06"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("E35").Value = "  -7.86%  "
$ws.Range("E36").Value = "  -4.55%  "
$ws.Range("D37").Value = "'7.00"
$ws.Range("E37").Value = "  +11.82%  "
$ws.Range("E38").Value = "  -11.27%  "
$ws.Range("D39").Value = "'22.31"
$ws.Range("E39").Value = "  -4.28%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "'21.66"
$ws.Range("E41").Value = "  +3.48%  "

# Rows 42/43 swap order: row 42 becomes PolygonEcosystemToken, row 43
# becomes Stacks (with refreshed price/volume data for each).
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.390"
$ws.Range("E42").Value = "  -5.30%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.98"
$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").Value = "'2.92"
$ws.Range("E44").Value = "  -4.73%  "
$ws.Range("D46").Value = "'155.19"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").Value = "'178.88"
$ws.Range("E47").Value = "  -6.43%  "
$ws.Range("D48").Value = "'44.79"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("E50").Value = "  -1.18%  "

# Row 51 changes from ARBITRUM to Stellar.
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.123"
$ws.Range("E51").Value = "  +11.66%  "
